$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.535309666666667
$ws.Range("H2").Value = 16.605929
$ws.Range("I2").Value = 0.1390876011186406
$ws.Range("J2").Value = 0.1461320332765693
$ws.Range("M2").Value = 51.02156433333334
$ws.Range("N2").Value = 153.064693
$ws.Range("O2").Value = 0.4760900215891154
$ws.Range("P2").Value = 0.4807539937572116
$ws.Range("Q2").Value = 282.4201582627553
$ws.Range("R2").Value = 2541.781424364797
$ws.Range("S2").Value = 0.06621821901935189
$ws.Range("T2").Value = 0.07025355861357242
$ws.Range("G3").Value = 5.535309666666667
$ws.Range("H3").Value = 16.605929
$ws.Range("I3").Value = 0.1390876011186406
$ws.Range("J3").Value = 0.1461320332765693
$ws.Range("O3").Value = 0.00890016878749362
$ws.Range("P3").Value = 0.008987358473548528
$ws.Range("Q3").Value = 5.279646628885888
$ws.Range("R3").Value = 47.516819659973
$ws.Range("S3").Value = 0.001237903126203488
$ws.Range("T3").Value = 0.00131334096752505
$ws.Range("G4").Value = 5.535309666666667
$ws.Range("H4").Value = 16.605929
$ws.Range("I4").Value = 0.1390876011186406
$ws.Range("J4").Value = 0.1461320332765693
$ws.Range("M4").Value = 28.11170133333333
$ws.Range("N4").Value = 84.335104
$ws.Range("O4").Value = 0.2623145854026591
$ws.Range("P4").Value = 0.2648843261452188
$ws.Range("Q4").Value = 155.6069721368462
$ws.Range("R4").Value = 1400.462749231616
$ws.Range("S4").Value = 0.03648470642208664
$ws.Range("T4").Value = 0.03870808516269474
$ws.Range("G5").Value = 5.535309666666667
$ws.Range("H5").Value = 16.605929
$ws.Range("I5").Value = 0.1390876011186406
$ws.Range("J5").Value = 0.1461320332765693
$ws.Range("M5").Value = 3.119026
$ws.Range("N5").Value = 6.238052
$ws.Range("O5").Value = 0.02910410872500189
$ws.Range("P5").Value = 0.01959281630196169
$ws.Range("Q5").Value = 17.26477476838467
$ws.Range("R5").Value = 103.588648610308
$ws.Range("S5").Value = 0.004048020665256611
$ws.Range("T5").Value = 0.002863138083819975
$ws.Range("G6").Value = 5.535309666666667
$ws.Range("H6").Value = 16.605929
$ws.Range("I6").Value = 0.1390876011186406
$ws.Range("J6").Value = 0.1461320332765693
$ws.Range("M6").Value = 23.96178866666667
$ws.Range("N6").Value = 71.885366
$ws.Range("O6").Value = 0.22359111549573
$ws.Range("P6").Value = 0.2257815053220593
$ws.Range("Q6").Value = 132.6359204372238
$ws.Range("R6").Value = 1193.723283935014
$ws.Range("S6").Value = 0.031098751885742
$ws.Range("T6").Value = 0.03299391044895707
$ws.Range("I7").Value = 0.2300921801028976
$ws.Range("J7").Value = 0.2417457620165159
$ws.Range("M7").Value = 51.02156433333334
$ws.Range("N7").Value = 153.064693
$ws.Range("O7").Value = 0.4760900215891154
$ws.Range("P7").Value = 0.4807539937572116
$ws.Range("Q7").Value = 467.2067775779165
$ws.Range("R7").Value = 4204.860998201249
$ws.Range("S7").Value = 0.1095445909926751
$ws.Range("T7").Value = 0.1162202405633205
$ws.Range("I8").Value = 0.2300921801028976
$ws.Range("J8").Value = 0.2417457620165159
$ws.Range("O8").Value = 0.00890016878749362
$ws.Range("P8").Value = 0.008987358473548528
$ws.Range("S8").Value = 0.002047859239598169
$ws.Range("T8").Value = 0.00217265582270358
$ws.Range("I9").Value = 0.2300921801028976
$ws.Range("J9").Value = 0.2417457620165159
$ws.Range("M9").Value = 28.11170133333333
$ws.Range("N9").Value = 84.335104
$ws.Range("O9").Value = 0.2623145854026591
$ws.Range("P9").Value = 0.2648843261452188
$ws.Range("Q9").Value = 257.4201235064605
$ws.Range("R9").Value = 2316.781111558144
$ws.Range("S9").Value = 0.06035653482808553
$ws.Range("T9").Value = 0.06403466327020726
$ws.Range("I10").Value = 0.2300921801028976
$ws.Range("J10").Value = 0.2417457620165159
$ws.Range("M10").Value = 3.119026
$ws.Range("N10").Value = 6.238052
$ws.Range("O10").Value = 0.02910410872500189
$ws.Range("P10").Value = 0.01959281630196169
$ws.Range("Q10").Value = 28.56106247784534
$ws.Range("R10").Value = 171.366374867072
$ws.Range("S10").Value = 0.006696627826487447
$ws.Range("T10").Value = 0.004736480306967344
$ws.Range("I11").Value = 0.2300921801028976
$ws.Range("J11").Value = 0.2417457620165159
$ws.Range("M11").Value = 23.96178866666667
$ws.Range("N11").Value = 71.885366
$ws.Range("O11").Value = 0.22359111549573
$ws.Range("P11").Value = 0.2257815053220593
$ws.Range("Q11").Value = 219.4191850884196
$ws.Range("R11").Value = 1974.772665795776
$ws.Range("S11").Value = 0.05144656721605128
$ws.Range("T11").Value = 0.05458172205331727
$ws.Range("G12").Value = 9.356602
$ws.Range("H12").Value = 28.069806
$ws.Range("I12").Value = 0.2351065080674274
$ws.Range("J12").Value = 0.2470140528999518
$ws.Range("M12").Value = 51.02156433333334
$ws.Range("N12").Value = 153.064693
$ws.Range("O12").Value = 0.4760900215891154
$ws.Range("P12").Value = 0.4807539937572116
$ws.Range("Q12").Value = 477.3884708843954
$ws.Range("R12").Value = 4296.496237959558
$ws.Range("S12").Value = 0.111931862501563
$ws.Range("T12").Value = 0.118752992445807
$ws.Range("G13").Value = 9.356602
$ws.Range("H13").Value = 28.069806
$ws.Range("I13").Value = 0.2351065080674274
$ws.Range("J13").Value = 0.2470140528999518
$ws.Range("O13").Value = 0.00890016878749362
$ws.Range("P13").Value = 0.008987358473548528
$ws.Range("Q13").Value = 8.924442385691332
$ws.Range("R13").Value = 80.319981471222
$ws.Range("S13").Value = 0.002092487604838334
$ws.Range("T13").Value = 0.002220003841415946
$ws.Range("G14").Value = 9.356602
$ws.Range("H14").Value = 28.069806
$ws.Range("I14").Value = 0.2351065080674274
$ws.Range("J14").Value = 0.2470140528999518
$ws.Range("M14").Value = 28.11170133333333
$ws.Range("N14").Value = 84.335104
$ws.Range("O14").Value = 0.2623145854026591
$ws.Range("P14").Value = 0.2648843261452188
$ws.Range("Q14").Value = 263.0300009188693
$ws.Range("R14").Value = 2367.270008269824
$ws.Range("S14").Value = 0.06167186618917413
$ws.Range("T14").Value = 0.06543015095080318
$ws.Range("G15").Value = 9.356602
$ws.Range("H15").Value = 28.069806
$ws.Range("I15").Value = 0.2351065080674274
$ws.Range("J15").Value = 0.2470140528999518
$ws.Range("M15").Value = 3.119026
$ws.Range("N15").Value = 6.238052
$ws.Range("O15").Value = 0.02910410872500189
$ws.Range("P15").Value = 0.01959281630196169
$ws.Range("Q15").Value = 29.183484909652
$ws.Range("R15").Value = 175.100909457912
$ws.Range("S15").Value = 0.00684256537274994
$ws.Range("T15").Value = 0.004839700962471803
$ws.Range("G16").Value = 9.356602
$ws.Range("H16").Value = 28.069806
$ws.Range("I16").Value = 0.2351065080674274
$ws.Range("J16").Value = 0.2470140528999518
$ws.Range("M16").Value = 23.96178866666667
$ws.Range("N16").Value = 71.885366
$ws.Range("O16").Value = 0.22359111549573
$ws.Range("P16").Value = 0.2257815053220593
$ws.Range("Q16").Value = 224.2009197621107
$ws.Range("R16").Value = 2017.808277858996
$ws.Range("S16").Value = 0.05256772639910193
$ws.Range("T16").Value = 0.05577120469945392
$ws.Range("G17").Value = 5.7553975
$ws.Range("H17").Value = 11.510795
$ws.Range("I17").Value = 0.144617822663078
$ws.Range("J17").Value = 0.1012948976223954
$ws.Range("M17").Value = 51.02156433333334
$ws.Range("N17").Value = 153.064693
$ws.Range("O17").Value = 0.4760900215891154
$ws.Range("P17").Value = 0.4807539937572116
$ws.Range("Q17").Value = 293.6493838101559
$ws.Range("R17").Value = 1761.896302860935
$ws.Range("S17").Value = 0.06885110231383568
$ws.Range("T17").Value = 0.04869792657919448
$ws.Range("G18").Value = 5.7553975
$ws.Range("H18").Value = 11.510795
$ws.Range("I18").Value = 0.144617822663078
$ws.Range("J18").Value = 0.1012948976223954
$ws.Range("O18").Value = 0.00890016878749362
$ws.Range("P18").Value = 0.008987358473548528
$ws.Range("Q18").Value = 5.489569118735832
$ws.Range("R18").Value = 32.93741471241499
$ws.Range("S18").Value = 0.001287123031381215
$ws.Range("T18").Value = 0.0009103735564738661
$ws.Range("G19").Value = 5.7553975
$ws.Range("H19").Value = 11.510795
$ws.Range("I19").Value = 0.144617822663078
$ws.Range("J19").Value = 0.1012948976223954
$ws.Range("M19").Value = 28.11170133333333
$ws.Range("N19").Value = 84.335104
$ws.Range("O19").Value = 0.2623145854026591
$ws.Range("P19").Value = 0.2648843261452188
$ws.Range("Q19").Value = 161.7940155746133
$ws.Range("R19").Value = 970.76409344768
$ws.Range("S19").Value = 0.03793536419370058
$ws.Range("T19").Value = 0.02683143069865714
$ws.Range("G20").Value = 5.7553975
$ws.Range("H20").Value = 11.510795
$ws.Range("I20").Value = 0.144617822663078
$ws.Range("J20").Value = 0.1012948976223954
$ws.Range("M20").Value = 3.119026
$ws.Range("N20").Value = 6.238052
$ws.Range("O20").Value = 0.02910410872500189
$ws.Range("P20").Value = 0.01959281630196169
$ws.Range("Q20").Value = 17.951234442835
$ws.Range("R20").Value = 71.80493777134
$ws.Range("S20").Value = 0.004208972834359265
$ws.Range("T20").Value = 0.001984652321441609
$ws.Range("G21").Value = 5.7553975
$ws.Range("H21").Value = 11.510795
$ws.Range("I21").Value = 0.144617822663078
$ws.Range("J21").Value = 0.1012948976223954
$ws.Range("M21").Value = 23.96178866666667
$ws.Range("N21").Value = 71.885366
$ws.Range("O21").Value = 0.22359111549573
$ws.Range("P21").Value = 0.2257815053220593
$ws.Range("Q21").Value = 137.9096185876617
$ws.Range("R21").Value = 827.45771152597
$ws.Range("S21").Value = 0.03233526028980128
$ws.Range("T21").Value = 0.02287051446662832
$ws.Range("G22").Value = 9.992936
$ws.Range("H22").Value = 29.978808
$ws.Range("I22").Value = 0.2510958880479564
$ws.Range("J22").Value = 0.2638132541845675
$ws.Range("M22").Value = 51.02156433333334
$ws.Range("N22").Value = 153.064693
$ws.Range("O22").Value = 0.4760900215891154
$ws.Range("P22").Value = 0.4807539937572116
$ws.Range("Q22").Value = 509.8552270028827
$ws.Range("R22").Value = 4588.697043025944
$ws.Range("S22").Value = 0.1195442467616897
$ws.Range("T22").Value = 0.1268292755553173
$ws.Range("G23").Value = 9.992936
$ws.Range("H23").Value = 29.978808
$ws.Range("I23").Value = 0.2510958880479564
$ws.Range("J23").Value = 0.2638132541845675
$ws.Range("O23").Value = 0.00890016878749362
$ws.Range("P23").Value = 0.008987358473548528
$ws.Range("Q23").Value = 9.531385603010666
$ws.Range("R23").Value = 85.78247042709599
$ws.Range("S23").Value = 0.002234795785472414
$ws.Range("T23").Value = 0.002370984285430085
$ws.Range("G24").Value = 9.992936
$ws.Range("H24").Value = 29.978808
$ws.Range("I24").Value = 0.2510958880479564
$ws.Range("J24").Value = 0.2638132541845675
$ws.Range("M24").Value = 28.11170133333333
$ws.Range("N24").Value = 84.335104
$ws.Range("O24").Value = 0.2623145854026591
$ws.Range("P24").Value = 0.2648843261452188
$ws.Range("Q24").Value = 280.9184322751146
$ws.Range("R24").Value = 2528.265890476032
$ws.Range("S24").Value = 0.06586611376961218
$ws.Range("T24").Value = 0.0698799960628565
$ws.Range("G25").Value = 9.992936
$ws.Range("H25").Value = 29.978808
$ws.Range("I25").Value = 0.2510958880479564
$ws.Range("J25").Value = 0.2638132541845675
$ws.Range("M25").Value = 3.119026
$ws.Range("N25").Value = 6.238052
$ws.Range("O25").Value = 0.02910410872500189
$ws.Range("P25").Value = 0.01959281630196169
$ws.Range("Q25").Value = 31.168227200336
$ws.Range("R25").Value = 187.009363202016
$ws.Range("S25").Value = 0.007307922026148626
$ws.Range("T25").Value = 0.005168844627260958
$ws.Range("G26").Value = 9.992936
$ws.Range("H26").Value = 29.978808
$ws.Range("I26").Value = 0.2510958880479564
$ws.Range("J26").Value = 0.2638132541845675
$ws.Range("M26").Value = 23.96178866666667
$ws.Range("N26").Value = 71.885366
$ws.Range("O26").Value = 0.22359111549573
$ws.Range("P26").Value = 0.2257815053220593
$ws.Range("Q26").Value = 239.4486205915254
$ws.Range("R26").Value = 2155.037585323728
$ws.Range("S26").Value = 0.05614280970503351
$ws.Range("T26").Value = 0.05956415365370272
